$d = $word.ActiveDocument

# Change 1: merge "5. " and the question into a single run (text unaffected,
# but this removes the split and the _GoBack bookmark wrapped around the
# question text).
$d.Content.Find.Execute("5. Qual a importância", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "5. Qual a importância", 2) | Out-Null

# Change 2: "indisponibilizar" -> "disponibilizar"
$d.Content.Find.Execute("indisponibilizar", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "disponibilizar", 2) | Out-Null
